# New sampling data added to several logger sheets.
$wb = $excel.ActiveWorkbook

# --- MO: bulk paste of 6 new logger-dump rows (rows 22-27) ---
$moWs = $wb.Worksheets.Item("MO")
$moWs.Activate() | Out-Null

$moRows = @(
    @(44296.304861111108, 0,  412.5, 8.6),
    @(44296.304166666669, 4,  412.9, 8.6),
    @(44296.302083333336, 8,  408.5, 8.1999999999999993),
    @(44296.299305555556, 12, 395.1, 7.1),
    @(44296.296527777777, 16, 390.6, 6.6),
    @(44296.293749999997, 20, 386.5, 6.6)
)

$r = 22
foreach ($row in $moRows) {
    $moWs.Cells.Item($r, 1).Value = $row[0]
    $moWs.Cells.Item($r, 2).Value = $row[1]
    $moWs.Cells.Item($r, 3).Value = $row[2]
    $moWs.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# --- WIC: one new row ---
$ws = $wb.Worksheets.Item("WIC")
$ws.Activate() | Out-Null
$ws.Cells.Item(28, 1).Value = 44296.418749999997
$ws.Cells.Item(28, 2).Value = 779.3
$ws.Cells.Item(28, 3).Value = 14
$ws.Range("F28").Select() | Out-Null

# --- YS: one new row ---
$ws = $wb.Worksheets.Item("YS")
$ws.Activate() | Out-Null
$ws.Cells.Item(43, 1).Value = 44296.442361111112
$ws.Cells.Item(43, 2).Value = 429
$ws.Cells.Item(43, 3).Value = 10.3
$ws.Range("C44").Select() | Out-Null

# --- SW: one new row ---
$ws = $wb.Worksheets.Item("SW")
$ws.Activate() | Out-Null
$ws.Cells.Item(40, 1).Value = 44296.474999999999
$ws.Cells.Item(40, 2).Value = 871.2
$ws.Cells.Item(40, 3).Value = 11.7
$ws.Range("C41").Select() | Out-Null

# --- YI: one new row ---
$ws = $wb.Worksheets.Item("YI")
$ws.Activate() | Out-Null
$ws.Cells.Item(40, 1).Value = 44296.490972222222
$ws.Cells.Item(40, 2).Value = 379.5
$ws.Cells.Item(40, 3).Value = 8.3000000000000007
$ws.Range("D42").Select() | Out-Null

# --- YN: one new row ---
$ws = $wb.Worksheets.Item("YN")
$ws.Activate() | Out-Null
$ws.Cells.Item(39, 1).Value = 44296.509722222225
$ws.Cells.Item(39, 2).Value = 586.6
$ws.Cells.Item(39, 3).Value = 14.4
$ws.Range("F34").Select() | Out-Null

# --- 6MC: one new row ---
$ws = $wb.Worksheets.Item("6MC")
$ws.Activate() | Out-Null
$ws.Cells.Item(41, 1).Value = 44296.520833333336
$ws.Cells.Item(41, 2).Value = 581
$ws.Cells.Item(41, 3).Value = 12.1
$ws.Range("C42").Select() | Out-Null

# --- DC: one new row ---
$ws = $wb.Worksheets.Item("DC")
$ws.Activate() | Out-Null
$ws.Cells.Item(41, 1).Value = 44296.527777777781
$ws.Cells.Item(41, 2).Value = 588.20000000000005
$ws.Cells.Item(41, 3).Value = 11.4
$ws.Range("C42").Select() | Out-Null

# --- PBMS: one new row ---
$ws = $wb.Worksheets.Item("PBMS")
$ws.Activate() | Out-Null
$ws.Cells.Item(42, 1).Value = 44296.540277777778
$ws.Cells.Item(42, 2).Value = 968.9
$ws.Cells.Item(42, 3).Value = 13.3
$ws.Range("C43").Select() | Out-Null

# --- PBSF: one new row ---
$ws = $wb.Worksheets.Item("PBSF")
$ws.Activate() | Out-Null
$ws.Cells.Item(42, 1).Value = 44296.552083333336
$ws.Cells.Item(42, 2).Value = 1145
$ws.Cells.Item(42, 3).Value = 13.8
$ws.Range("G42").Select() | Out-Null

# --- Session ends back on MO, cell H6 selected ---
$moWs.Activate() | Out-Null
$moWs.Range("H6").Select() | Out-Null
